$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "261.02"
$c.Style = "Normal"

$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value = "-0.04%"
$c.Style = "Normal"

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "26.84"
$c.Style = "Normal"

$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value = "-1.72%"
$c.Style = "Normal"

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "4.711"
$c.Style = "Normal"

$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = "@"
$c.Value = "0.15%"
$c.Style = "Normal"

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "0.06221"
$c.Style = "Normal"

$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value = "2.33%"
$c.Style = "Normal"

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "6.733"
$c.Style = "Normal"

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.8504"
$c.Style = "Normal"

$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value = "0.56%"
$c.Style = "Normal"

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.9109"
$c.Style = "Normal"

$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value = "-1.37%"
$c.Style = "Normal"

$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value = "-0.22%"
$c.Style = "Normal"

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "0.04959"
$c.Style = "Normal"

$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value = "-1.53%"
$c.Style = "Normal"

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.07076"
$c.Style = "Normal"

$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value = "-0.40%"
$c.Style = "Normal"

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "0.03083"
$c.Style = "Normal"

$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value = "-1.51%"
$c.Style = "Normal"

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "0.09052"
$c.Style = "Normal"

$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = "@"
$c.Value = "-0.22%"
$c.Style = "Normal"

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "0.001538"
$c.Style = "Normal"

$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value = "0.10%"
$c.Style = "Normal"

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "0.0006173"
$c.Style = "Normal"

$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value = "1.68%"
$c.Style = "Normal"

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "0.005966"
$c.Style = "Normal"

$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = "@"
$c.Value = "-2.53%"
$c.Style = "Normal"

$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = "@"
$c.Value = "-0.13%"
$c.Style = "Normal"

$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value = "0.80%"
$c.Style = "Normal"

$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = "@"
$c.Value = "0.03%"
$c.Style = "Normal"

$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = "@"
$c.Value = "-0.65%"
$c.Style = "Normal"

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "0.1310"
$c.Style = "Normal"

$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = "@"
$c.Value = "0.32%"
$c.Style = "Normal"

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "4.120"
$c.Style = "Normal"

$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = "@"
$c.Value = "0.73%"
$c.Style = "Normal"

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "0.04241"
$c.Style = "Normal"

$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = "@"
$c.Value = "0.16%"
$c.Style = "Normal"

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "0.001201"
$c.Style = "Normal"

$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = "@"
$c.Value = "-1.55%"
$c.Style = "Normal"

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "0.004071"
$c.Style = "Normal"

$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = "@"
$c.Value = "4.09%"
$c.Style = "Normal"

$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = "@"
$c.Value = "0.02%"
$c.Style = "Normal"

$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = "@"
$c.Value = "4.09%"
$c.Style = "Normal"

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "0.03936"
$c.Style = "Normal"

$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = "@"
$c.Value = "1.70%"
$c.Style = "Normal"

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "0.1112"
$c.Style = "Normal"

$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = "@"
$c.Value = "-0.19%"
$c.Style = "Normal"

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.004126"
$c.Style = "Normal"

$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value = "0.75%"
$c.Style = "Normal"

$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value = "-5.40%"
$c.Style = "Normal"

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.01337"
$c.Style = "Normal"

$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = "@"
$c.Value = "-18.25%"
$c.Style = "Normal"

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "0.00005163"
$c.Style = "Normal"

$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = "@"
$c.Value = "-2.15%"
$c.Style = "Normal"

$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = "@"
$c.Value = "0.01%"
$c.Style = "Normal"

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "0.2479"
$c.Style = "Normal"

$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value = "83.20%"
$c.Style = "Normal"

$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = "@"
$c.Value = "0.01%"
$c.Style = "Normal"

$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value = "0.01%"
$c.Style = "Normal"
